$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("survey")
$ws2 = $wb.Worksheets.Item("properties")
$ws3 = $wb.Worksheets.Item("settings")

# Rename the translation keys (content-level change: "display.text",
# "display.hint", "display.title" -> "*.text" suffixed variants).
# Update settings!C1 ("display.title" -> "display.title.text") first so the
# new shared strings land in the same order as the target workbook.
$ws3.Range("C2").Select()
$ws3.Range("C1").Value = "display.title.text"

$ws1.Range("H1").Value = "display.prompt.text"
$ws1.Range("I1").Value = "display.hint.text"

# Move the active sheet / selection to match the new view state: "survey"
# becomes the selected tab with B10 highlighted, "properties" loses the
# tab-selected flag (its own selection, E5, is untouched).
$ws1.Activate()
$ws1.Range("B10").Select()
